# binance_ETHUSDT_data.xlsx — row 20 update (Sheet1)
# Text-valued columns keep their trailing-zero string formatting, so they are
# entered with a leading apostrophe (quote-prefix) to force literal text,
# exactly as Excel itself does when a number-looking value must stay text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# close
$ws.Range("E20").Value = "'127.12000000"
# volume
$ws.Range("F20").Value = "'84693.74392000"
# Quote asset volume
$ws.Range("H20").Value = "'10788650.59353520"
# number of trades
$ws.Range("I20").Value = 37982
# Taker buy base asset volume
$ws.Range("J20").Value = "'43316.73631000"
# Taker buy quote asset volume
$ws.Range("K20").Value = "'5519382.19591250"

# MA_1
$ws.Range("M20").Value = 127.1199999999999
# MA_2
$ws.Range("N20").Value = 126.7050000000001
# MA_7
$ws.Range("O20").Value = 127.3885714285714
# MA_15
$ws.Range("P20").Value = 129.776

# EMA_1
$ws.Range("R20").Value = 127.12
# EMA_2
$ws.Range("S20").Value = 126.78311699612
# EMA_12
$ws.Range("T20").Value = 128.4922871171791
# EMA_26
$ws.Range("U20").Value = 130.4321812872458
# DIF
$ws.Range("V20").Value = -1.939894170066708
# DEM
$ws.Range("W20").Value = -1.811387779898418
# OSC
$ws.Range("X20").Value = -0.1285063901682899
